$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") updates. The source data stores these as plain text
# cells (t="inlineStr"), not numbers, so we force text interpretation via
# NumberFormat "@" before assigning each Value (otherwise Excel
# auto-converts the numeric-looking string to a real number and we lose
# the exact text representation, e.g. trailing zeros like "0.08330").
# The format is reset back to the default "Normal" style afterwards so
# each cell keeps its original (unstyled) appearance.
$updates = @(
    @{ Row = 2; Value = "281.07" },
    @{ Row = 3; Value = "20.78" },
    @{ Row = 4; Value = "6.262" },
    @{ Row = 5; Value = "0.06139" },
    @{ Row = 6; Value = "3.575" },
    @{ Row = 7; Value = "6.563" },
    @{ Row = 8; Value = "1.478" },
    @{ Row = 9; Value = "0.8177" },
    @{ Row = 10; Value = "0.01377" },
    @{ Row = 11; Value = "0.1625" },
    @{ Row = 12; Value = "0.08323" },
    @{ Row = 13; Value = "0.03530" },
    @{ Row = 14; Value = "0.03208" },
    @{ Row = 15; Value = "0.09140" },
    @{ Row = 16; Value = "3.712" },
    @{ Row = 17; Value = "0.001644" },
    @{ Row = 18; Value = "0.04644" },
    @{ Row = 19; Value = "0.006424" },
    @{ Row = 20; Value = "0.006156" },
    @{ Row = 22; Value = "0.0001502" },
    @{ Row = 23; Value = "3.801" },
    @{ Row = 24; Value = "2.294" },
    @{ Row = 25; Value = "0.3335" },
    @{ Row = 26; Value = "0.1251" },
    @{ Row = 40; Value = "0.04674" },
    @{ Row = 41; Value = "0.006310" },
    @{ Row = 42; Value = "0.007179" },
    @{ Row = 43; Value = "0.1097" },
    @{ Row = 44; Value = "0.01133" },
    @{ Row = 45; Value = "0.00006365" },
    @{ Row = 46; Value = "0.00000000751" },
    @{ Row = 47; Value = "1.002" },
    @{ Row = 48; Value = "0.002946" },
    @{ Row = 49; Value = "0.00001903" },
    @{ Row = 50; Value = "0.01242" }
)

foreach ($u in $updates) {
    $cell = $ws.Range("D" + $u.Row)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
